$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '93.247.57'
$ws.Range("E2").Value = '  +1.90%  '

$ws.Range("D3").Value = '3.124.58'
$ws.Range("E3").Value = '  +0.44%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").Value = '''237.91'
$ws.Range("E5").Value = '  -2.25%  '

$ws.Range("D6").Value = '''613.95'
$ws.Range("E6").Value = '  -0.31%  '

$ws.Range("D7").Value = '''1.10'
$ws.Range("E7").Value = '  +1.01%  '

$ws.Range("E8").Value = '  +1.86%  '

$ws.Range("D9").Value = '''1.00'
$ws.Range("E9").Value = '  +0.01%  '

$ws.Range("D10").Value = '''0.839'
$ws.Range("E10").Value = '  +13.36%  '

$ws.Range("D11").Value = '3.124.28'
$ws.Range("E11").Value = '  +0.50%  '

$ws.Range("E12").Value = '  -2.38%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").Value = '''35.44'
$ws.Range("E13").Value = '  +2.74%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").Value = '''0.0000246'
$ws.Range("E14").Value = '  -1.26%  '

$ws.Range("D15").Value = '93.071.42'
$ws.Range("E15").Value = '  +1.90%  '

$ws.Range("D16").Value = '''5.45'
$ws.Range("E16").Value = '  -2.92%  '

$ws.Range("D17").Value = '3.710.30'
$ws.Range("E17").Value = '  +0.73%  '

$ws.Range("D18").Value = '3.125.03'
$ws.Range("E18").Value = '  +0.42%  '

$ws.Range("E19").Value = '  +1.08%  '

$ws.Range("D20").Value = '''14.84'
$ws.Range("E20").Value = '  +0.78%  '

$ws.Range("E21").Value = '  +4.91%  '

$ws.Range("E22").Value = '  +0.60%  '

$ws.Range("D23").Value = '''443.17'
$ws.Range("E23").Value = '  -0.10%  '

$ws.Range("D24").Value = '''9.17'
$ws.Range("E24").Value = '  -1.41%  '

$ws.Range("B25").Value = 'Aptos'
$ws.Range("C25").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D25").Value = '''13.16'
$ws.Range("E25").Value = '  +12.88%  '

$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D26").Value = '''5.75'
$ws.Range("E26").Value = '  +2.24%  '

$ws.Range("D27").Value = '''86.17'

$ws.Range("B28").Value = 'Dai'
$ws.Range("C28").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D28").Value = '''0.999'
$ws.Range("E28").Value = '  +0.31%  '

$ws.Range("B29").Value = 'Cronos'
$ws.Range("C29").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D29").Value = '''0.182'
$ws.Range("E29").Value = '  +10.01%  '

$ws.Range("B30").Value = 'Stellar'
$ws.Range("C30").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D30").Value = '''0.239'
$ws.Range("E30").Value = '  +3.98%  '

$ws.Range("B31").Value = 'Hedera'
$ws.Range("C31").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D31").Value = '''0.123'
$ws.Range("E31").Value = '  -13.20%  '

$ws.Range("B32").Value = 'InternetComputer(DFINITY)'
$ws.Range("C32").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D32").Value = '''9.30'
$ws.Range("E32").Value = '  +0.22%  '

$ws.Range("E33").Value = '  +4.81%  '

$ws.Range("B34").Value = 'RenderToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D34").Value = '''8.09'
$ws.Range("E34").Value = '  +6.01%  '

$ws.Range("B35").Value = 'Kaspa'
$ws.Range("C35").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D35").Value = '''0.160'
$ws.Range("E35").Value = '  -9.35%  '

$ws.Range("B36").Value = 'EthereumClassic'
$ws.Range("C36").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D36").Value = '''26.12'
$ws.Range("E36").Value = '  -0.05%  '

$ws.Range("B37").Value = 'MantraDAO'
$ws.Range("C37").Value = 'https://coinranking.com/coin/cTdD8lD-6+mantradao-om'
$ws.Range("D37").Value = '''3.99'
$ws.Range("E37").Value = '  -1.62%  '

$ws.Range("B38").Value = 'PancakeSwap'
$ws.Range("C38").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D38").Value = '''1.91'
$ws.Range("E38").Value = '  -0.90%  '

$ws.Range("B39").Value = 'PolygonEcosystemToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
$ws.Range("D39").Value = '''0.453'
$ws.Range("E39").Value = '  +4.20%  '

$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").Value = '''1.30'
$ws.Range("E40").Value = '  +0.33%  '

$ws.Range("B41").Value = 'Bittensor'
$ws.Range("C41").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D41").Value = '''477.29'
$ws.Range("E41").Value = '  -2.34%  '

$ws.Range("E42").Value = '  +8.06%  '

$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '''3.34'
$ws.Range("E43").Value = '  -1.97%  '

$ws.Range("B44").Value = 'USDe'
$ws.Range("C44").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D44").Value = '''1.00'
$ws.Range("E44").Value = '  +0.03%  '

$ws.Range("B45").Value = 'Monero'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D45").Value = '''158.85'
$ws.Range("E45").Value = '  -0.30%  '

$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''0.703'
$ws.Range("E46").Value = '  +1.03%  '

$ws.Range("B47").Value = 'Stacks'
$ws.Range("C47").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D47").Value = '''1.86'
$ws.Range("E47").Value = '  -1.69%  '

$ws.Range("B48").Value = 'ImmutableX'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D48").Value = '''1.34'
$ws.Range("E48").Value = '  +0.66%  '

$ws.Range("B49").Value = 'Filecoin'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D49").Value = '''4.46'
$ws.Range("E49").Value = '  +1.88%  '

$ws.Range("B50").Value = 'OKB'
$ws.Range("C50").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D50").Value = '''44.07'
$ws.Range("E50").Value = '  +0.03%  '

$ws.Range("B51").Value = 'VeChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D51").Value = '''0.0324'
$ws.Range("E51").Value = '  +2.19%  '

